# Performance.xlsx - "Rearrange the repository, added entry for values of the RQ1"
#
# 1) Update header row labels (A1, B1) to the new wording.
# 2) Add a new data row (48) continuing the Paper/Pages/Time/Rate table, copying
#    the formatting of the last existing data row (46) so styles + the D-column
#    ratio formula match.
# 3) Add a further new row (49) with just Task/Pages/Time filled in (no Rate).
# 4) Move the active selection to B49 to match where the author ended up.
#
# NOTE: the order the new text values are written in matters, because it
# determines the order new entries are appended to the shared-strings table
# (xl/sharedStrings.xml) - we write them in the same order the target
# workbook has them in (Task and Paper, Fill up RQ1, Pages or Items,
# Work on the values of RQ1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting of the last populated row down onto the two new rows -
# Row 48 keeps the full 4-column (Paper/Pages/Time/Rate) pattern, but row 49
# only has Paper/Pages/Time filled in (no Rate computed yet), so only copy
# the first three columns' formatting there.
$ws.Range("A46:D46").Copy()
$ws.Range("A48:D48").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A46:C46").Copy()
$ws.Range("A49:C49").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- 1. Header row relabeling, and new row content, in shared-string order -
$ws.Range("A1").Value = "Task and Paper"
$ws.Range("A48").Value = "Fill up RQ1"
$ws.Range("B1").Value = "Pages or Items"
$ws.Range("A49").Value = "Work on the values of RQ1"

# --- 2. Row 48: Pages=20, Time=60, Rate=Time/Pages -------------------------
$ws.Range("B48").Value = 20
$ws.Range("C48").Value = 60
$ws.Range("D48").Formula = "=C48/B48"

# --- 3. Row 49: Pages=1, Time=90, no Rate entered yet ----------------------
$ws.Range("B49").Value = 1
$ws.Range("C49").Value = 90

# --- 4. Leave the selection where the author left it -----------------------
$ws.Range("B49").Select() | Out-Null
